# The document ends with a paragraph containing the "STATES" picture
# (which carries the Word-managed "_GoBack" bookmark marking the most
# recent edit location) followed by five empty trailing paragraphs and
# the section properties.
#
# The edit removes the very last (now superfluous) empty trailing
# paragraph and relocates the "_GoBack" bookmark from the picture
# paragraph down onto the new last empty paragraph - i.e. Word's
# auto-tracked "last edit" spot moved further down the document as
# more work (the mockup/backend additions mentioned in the commit)
# happened after the picture was inserted.

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$secondLastPara = $d.Paragraphs.Item($count - 1)

# Delete the final empty paragraph by removing its paragraph mark
# together with the preceding one, merging it away.
$deleteRange = $d.Range($secondLastPara.Range.End - 1, $lastPara.Range.End)
$deleteRange.Delete()

# Re-seat on the (now) last paragraph in the document and plant the
# "_GoBack" bookmark there; Word keeps only a single "_GoBack" bookmark,
# so adding it here automatically removes it from the picture paragraph.
$newLastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $newLastPara.Range)
